# aggiornamento fino a 20/09/2021
# Append new daily rows (375-385) to Sheet1, continuing the existing series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(375, 44449, 1, 7, 106.6098081023454),
    @(376, 44450, 0, 7, 106.6098081023454),
    @(377, 44451, 0, 6, 91.37983551629607),
    @(378, 44452, 0, 1, 15.22997258604935),
    @(379, 44453, 0, 1, 15.22997258604935),
    @(380, 44454, 0, 1, 15.22997258604935),
    @(381, 44455, 1, 2, 30.45994517209869),
    @(382, 44456, 0, 1, 15.22997258604935),
    @(383, 44457, 3, 4, 60.91989034419738),
    @(384, 44458, 1, 5, 76.14986293024673),
    @(385, 44459, 0, 5, 76.14986293024673)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    # Copy the formatting of the last existing data row (374) down onto the
    # new row so the date column keeps the same style (border/alignment/numfmt)
    $ws.Range("A374").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $rowData[1]
    $ws.Cells.Item($r, 2).Value = $rowData[2]
    $ws.Cells.Item($r, 3).Value = $rowData[3]
    $ws.Cells.Item($r, 4).Value = $rowData[4]
}
